$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 267.46667
$ws.Range("J17").Value = 267.46667
$ws.Range("L17").Value = 802.4000100000001
$ws.Range("N17").Value = -1138.40001
$ws.Range("H19").Value = 790
$ws.Range("I19").Value = 500.14285
$ws.Range("J19").Value = 1015.44446
$ws.Range("K19").Value = 500.14285
$ws.Range("L19").Value = 1015.44446
$ws.Range("M19").Value = -325.14285
$ws.Range("N19").Value = -1365.44446
$ws.Range("H62").Value = 66674348
$ws.Range("I62").Value = 71436580
$ws.Range("J62").Value = 3100
$ws.Range("K62").Value = 71436580
$ws.Range("L62").Value = 3100
$ws.Range("M62").Value = -71435956
$ws.Range("N62").Value = -4348
$ws.Range("H65").Value = 66674348
$ws.Range("I65").Value = 71436580
$ws.Range("J65").Value = 3100
$ws.Range("K65").Value = 357182900
$ws.Range("L65").Value = 15500
$ws.Range("M65").Value = -357179780
$ws.Range("N65").Value = -21740
$ws.Range("H82").Value = 2040.3334
$ws.Range("I82").Value = 2040.3334
$ws.Range("K82").Value = 6121.0002
$ws.Range("M82").Value = -5715.0002
$ws.Range("H85").Value = 2040.3334
$ws.Range("I85").Value = 2040.3334
$ws.Range("K85").Value = 6121.0002
$ws.Range("M85").Value = -4717.0002
$ws.Range("H98").Value = 16177.675
$ws.Range("I98").Value = 17769.781
$ws.Range("J98").Value = 9809.25
$ws.Range("K98").Value = 17769.781
$ws.Range("L98").Value = 9809.25
$ws.Range("M98").Value = -16271.781
$ws.Range("N98").Value = -12805.25
$ws.Range("H99").Value = 2346149.5
$ws.Range("I99").Value = 3280409.2
$ws.Range("J99").Value = 10500
$ws.Range("K99").Value = 9841227.600000001
$ws.Range("L99").Value = 31500
$ws.Range("M99").Value = -9839729.600000001
$ws.Range("N99").Value = -34496
$ws.Range("H101").Value = 9524955
$ws.Range("I101").Value = 10990110
$ws.Range("J101").Value = 1447
$ws.Range("K101").Value = 32970330
$ws.Range("L101").Value = 4341
$ws.Range("M101").Value = -32968708
$ws.Range("N101").Value = -7585
$ws.Range("H104").Value = 673.5
$ws.Range("I104").Value = 231.33333
$ws.Range("K104").Value = 693.99999
$ws.Range("M104").Value = 1053.00001
$ws.Range("H106").Value = 7265263
$ws.Range("I106").Value = 8820198
$ws.Range("J106").Value = 8900
$ws.Range("K106").Value = 8820198
$ws.Range("L106").Value = 8900
$ws.Range("M106").Value = -8819567
$ws.Range("N106").Value = -10162
$ws.Range("H107").Value = 8204.154
$ws.Range("I107").Value = 10445.8
$ws.Range("J107").Value = 732
$ws.Range("K107").Value = 10445.8
$ws.Range("L107").Value = 732
$ws.Range("M107").Value = -8525.799999999999
$ws.Range("N107").Value = -4572
$ws.Range("H115").Value = 557.5
$ws.Range("I115").Value = 225
$ws.Range("J115").Value = 668.3333
$ws.Range("K115").Value = 675
$ws.Range("L115").Value = 2004.9999
$ws.Range("M115").Value = 892
$ws.Range("N115").Value = -5138.9999
$ws.Range("H118").Value = 1160.125
$ws.Range("I118").Value = 1111.5714
$ws.Range("K118").Value = 3334.7142
$ws.Range("M118").Value = -1677.7142
$ws.Range("H122").Value = 16177.675
$ws.Range("I122").Value = 17769.781
$ws.Range("J122").Value = 9809.25
$ws.Range("K122").Value = 53309.34299999999
$ws.Range("L122").Value = 29427.75
$ws.Range("M122").Value = -50859.34299999999
$ws.Range("N122").Value = -34327.75
$ws.Range("H127").Value = 1742.375
$ws.Range("I127").Value = 1250.8334
$ws.Range("J127").Value = 3217
$ws.Range("K127").Value = 3752.5002
$ws.Range("L127").Value = 9651
$ws.Range("M127").Value = 1207.4998
$ws.Range("N127").Value = -19571
$ws.Range("H132").Value = 4891.364
$ws.Range("I132").Value = 4152
$ws.Range("K132").Value = 12456
$ws.Range("M132").Value = -9926
$ws.Range("H137").Value = 348614.94
$ws.Range("I137").Value = 543445.2
$ws.Range("J137").Value = 12089.909
$ws.Range("K137").Value = 1630335.6
$ws.Range("L137").Value = 36269.727
$ws.Range("M137").Value = -1627785.6
$ws.Range("N137").Value = -41369.727
$ws.Range("H138").Value = 7302.5713
$ws.Range("I138").Value = 3609
$ws.Range("J138").Value = 7918.1665
$ws.Range("K138").Value = 10827
$ws.Range("L138").Value = 23754.4995
$ws.Range("M138").Value = -5687
$ws.Range("N138").Value = -34034.49950000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H45").Value = 76043.97
$ws.Range("I45").Value = 94757.74000000001
$ws.Range("K45").Value = 94757.74000000001
$ws.Range("M45").Value = -94380.74000000001
$ws.Range("H61").Value = 11621.387
$ws.Range("I61").Value = 13779.826
$ws.Range("J61").Value = 5415.875
$ws.Range("K61").Value = 13779.826
$ws.Range("L61").Value = 5415.875
$ws.Range("M61").Value = -13567.826
$ws.Range("N61").Value = -5839.875
$ws.Range("H74").Value = 1319.0444
$ws.Range("I74").Value = 1205.7028
$ws.Range("J74").Value = 1843.25
$ws.Range("K74").Value = 1205.7028
$ws.Range("L74").Value = 1843.25
$ws.Range("M74").Value = -331.7028
$ws.Range("N74").Value = -3591.25
$ws.Range("H77").Value = 1319.0444
$ws.Range("I77").Value = 1205.7028
$ws.Range("J77").Value = 1843.25
$ws.Range("K77").Value = 6028.514
$ws.Range("L77").Value = 9216.25
$ws.Range("M77").Value = -1660.514
$ws.Range("N77").Value = -17952.25
$ws.Range("H125").Value = 170800
$ws.Range("J125").Value = 170800
$ws.Range("L125").Value = 170800
$ws.Range("N125").Value = -180640
$ws.Range("H136").Value = 11621.387
$ws.Range("I136").Value = 13779.826
$ws.Range("J136").Value = 5415.875
$ws.Range("K136").Value = 41339.478
$ws.Range("L136").Value = 16247.625
$ws.Range("M136").Value = -38789.478
$ws.Range("N136").Value = -21347.625

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 394.6
$ws.Range("I22").Value = 395.75
$ws.Range("J22").Value = 390
$ws.Range("K22").Value = 395.75
$ws.Range("L22").Value = 390
$ws.Range("M22").Value = -222.75
$ws.Range("N22").Value = -736
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H94").Value = 9445.293
$ws.Range("I94").Value = 11279.7
$ws.Range("K94").Value = 11279.7
$ws.Range("M94").Value = -10828.7
$ws.Range("H99").Value = 17309.564
$ws.Range("I99").Value = 17710
$ws.Range("K99").Value = 17710
$ws.Range("M99").Value = -16212
$ws.Range("H134").Value = 6725.2905
$ws.Range("I134").Value = 6844.276
$ws.Range("K134").Value = 20532.828
$ws.Range("M134").Value = -17997.828

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49780
$ws.Range("J20").Value = 49780
$ws.Range("L20").Value = 49780
$ws.Range("N20").Value = -50252
$ws.Range("H22").Value = 12820712
$ws.Range("I22").Value = 232.33333
$ws.Range("K22").Value = 232.33333
$ws.Range("M22").Value = 117.66667
$ws.Range("H30").Value = 49780
$ws.Range("J30").Value = 49780
$ws.Range("L30").Value = 49780
$ws.Range("N30").Value = -49962
$ws.Range("H31").Value = 3122.3865
$ws.Range("I31").Value = 2321.3333
$ws.Range("J31").Value = 3248.8684
$ws.Range("K31").Value = 2321.3333
$ws.Range("L31").Value = 3248.8684
$ws.Range("M31").Value = -2026.3333
$ws.Range("N31").Value = -3838.8684
$ws.Range("H34").Value = 3122.3865
$ws.Range("I34").Value = 2321.3333
$ws.Range("J34").Value = 3248.8684
$ws.Range("K34").Value = 2321.3333
$ws.Range("L34").Value = 3248.8684
$ws.Range("M34").Value = -2119.3333
$ws.Range("N34").Value = -3652.8684
$ws.Range("H99").Value = 211945.12
$ws.Range("I99").Value = 457840.62
$ws.Range("J99").Value = 3879.6924
$ws.Range("K99").Value = 457840.62
$ws.Range("L99").Value = 3879.6924
$ws.Range("M99").Value = -456342.62
$ws.Range("N99").Value = -6875.6924
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H114").Value = 54999.5
$ws.Range("J114").Value = 54999.5
$ws.Range("L114").Value = 54999.5
$ws.Range("N114").Value = -63677.5
$ws.Range("H115").Value = 54499.5
$ws.Range("J115").Value = 54499.5
$ws.Range("L115").Value = 54499.5
$ws.Range("N115").Value = -56849.5
$ws.Range("H120").Value = 50148
$ws.Range("I120").Value = 46962.668
$ws.Range("K120").Value = 46962.668
$ws.Range("M120").Value = -43333.668
$ws.Range("H121").Value = 21874
$ws.Range("J121").Value = 22663
$ws.Range("L121").Value = 22663
$ws.Range("N121").Value = -25283
$ws.Range("H126").Value = 211945.12
$ws.Range("I126").Value = 457840.62
$ws.Range("J126").Value = 3879.6924
$ws.Range("K126").Value = 1373521.86
$ws.Range("L126").Value = 11639.0772
$ws.Range("M126").Value = -1371051.86
$ws.Range("N126").Value = -16579.0772
$ws.Range("H128").Value = 49780
$ws.Range("J128").Value = 49780
$ws.Range("L128").Value = 49780
$ws.Range("N128").Value = -59740
$ws.Range("H134").Value = 3298783.5
$ws.Range("I134").Value = 20874238
$ws.Range("J134").Value = 3385.75
$ws.Range("K134").Value = 62622714
$ws.Range("L134").Value = 10157.25
$ws.Range("M134").Value = -62620179
$ws.Range("N134").Value = -15227.25
$ws.Range("H137").Value = 99995
$ws.Range("J137").Value = 99995
$ws.Range("L137").Value = 99995
$ws.Range("N137").Value = -110195

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 33665
$ws.Range("J52").Value = 33665
$ws.Range("L52").Value = 100995
$ws.Range("N52").Value = -101527
$ws.Range("H68").Value = 4135.6665
$ws.Range("I68").Value = 1981
$ws.Range("J68").Value = 4566.6
$ws.Range("K68").Value = 5943
$ws.Range("L68").Value = 13699.8
$ws.Range("M68").Value = -5132
$ws.Range("N68").Value = -15321.8
$ws.Range("H71").Value = 4135.6665
$ws.Range("I71").Value = 1981
$ws.Range("J71").Value = 4566.6
$ws.Range("K71").Value = 17829
$ws.Range("L71").Value = 41099.4
$ws.Range("M71").Value = -13773
$ws.Range("N71").Value = -49211.4
$ws.Range("H107").Value = 2271.3845
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 2385.6667
$ws.Range("K107").Value = 2700
$ws.Range("L107").Value = 7157.000100000001
$ws.Range("M107").Value = -780
$ws.Range("N107").Value = -10997.0001
$ws.Range("H113").Value = 2392.5833
$ws.Range("J113").Value = 2492.0908
$ws.Range("L113").Value = 7476.2724
$ws.Range("N113").Value = -11816.2724
$ws.Range("H122").Value = 1956.5
$ws.Range("I122").Value = 825.6
$ws.Range("J122").Value = 2360.3928
$ws.Range("K122").Value = 7430.400000000001
$ws.Range("L122").Value = 21243.5352
$ws.Range("M122").Value = -4980.400000000001
$ws.Range("N122").Value = -26143.5352
$ws.Range("H131").Value = 5172.857
$ws.Range("J131").Value = 1815.5883
$ws.Range("L131").Value = 5446.7649
$ws.Range("N131").Value = -15526.7649
$ws.Range("H132").Value = 13931997
$ws.Range("J132").Value = 15198456
$ws.Range("L132").Value = 136786104
$ws.Range("N132").Value = -136791164

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 198790.92
$ws.Range("J20").Value = 7207.4546
$ws.Range("L20").Value = 7207.4546
$ws.Range("N20").Value = -7697.4546
$ws.Range("H24").Value = 806405.6
$ws.Range("J24").Value = 8007
$ws.Range("L24").Value = 8007
$ws.Range("N24").Value = -8353
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H43").Value = 6654.857
$ws.Range("I43").Value = 6654.857
$ws.Range("K43").Value = 6654.857
$ws.Range("M43").Value = -6503.857
$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31136
$ws.Range("H102").Value = 26161.166
$ws.Range("I102").Value = 30793.4
$ws.Range("K102").Value = 30793.4
$ws.Range("M102").Value = -29171.4
$ws.Range("H122").Value = 6412.433
$ws.Range("I122").Value = 6588.6787
$ws.Range("K122").Value = 19766.0361
$ws.Range("M122").Value = -17316.0361
$ws.Range("H132").Value = 5155.4287
$ws.Range("I132").Value = 4652.7915
$ws.Range("K132").Value = 13958.3745
$ws.Range("M132").Value = -11428.3745

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 51273.223
$ws.Range("I7").Value = 63426.285
$ws.Range("K7").Value = 63426.285
$ws.Range("M7").Value = -63314.285
$ws.Range("H22").Value = 10825.462
$ws.Range("I22").Value = 14593.462
$ws.Range("J22").Value = 7057.4614
$ws.Range("K22").Value = 14593.462
$ws.Range("L22").Value = 7057.4614
$ws.Range("M22").Value = -14298.462
$ws.Range("N22").Value = -7647.4614
$ws.Range("H27").Value = 10825.462
$ws.Range("I27").Value = 14593.462
$ws.Range("J27").Value = 7057.4614
$ws.Range("K27").Value = 14593.462
$ws.Range("L27").Value = 7057.4614
$ws.Range("M27").Value = -14486.462
$ws.Range("N27").Value = -7271.4614
$ws.Range("H68").Value = 4397.643
$ws.Range("I68").Value = 2228
$ws.Range("J68").Value = 6024.875
$ws.Range("K68").Value = 2228
$ws.Range("L68").Value = 6024.875
$ws.Range("M68").Value = -1479
$ws.Range("N68").Value = -7522.875
$ws.Range("H71").Value = 4397.643
$ws.Range("I71").Value = 2228
$ws.Range("J71").Value = 6024.875
$ws.Range("K71").Value = 11140
$ws.Range("L71").Value = 30124.375
$ws.Range("M71").Value = -7396
$ws.Range("N71").Value = -37612.375
$ws.Range("H82").Value = 3278.1052
$ws.Range("I82").Value = 3357.9092
$ws.Range("J82").Value = 3168.375
$ws.Range("K82").Value = 3357.9092
$ws.Range("L82").Value = 3168.375
$ws.Range("M82").Value = -2996.9092
$ws.Range("N82").Value = -3890.375
$ws.Range("H85").Value = 3278.1052
$ws.Range("I85").Value = 3357.9092
$ws.Range("J85").Value = 3168.375
$ws.Range("K85").Value = 3357.9092
$ws.Range("L85").Value = 3168.375
$ws.Range("M85").Value = -2109.9092
$ws.Range("N85").Value = -5664.375
$ws.Range("H93").Value = 2863
$ws.Range("I93").Value = 2861.3914
$ws.Range("K93").Value = 2861.3914
$ws.Range("M93").Value = -1613.3914
$ws.Range("H122").Value = 5204.875
$ws.Range("I122").Value = 4597
$ws.Range("K122").Value = 13791
$ws.Range("M122").Value = -11341
$ws.Range("H126").Value = 51273.223
$ws.Range("I126").Value = 63426.285
$ws.Range("K126").Value = 190278.855
$ws.Range("M126").Value = -187808.855
$ws.Range("H132").Value = 11410.161
$ws.Range("I132").Value = 16088.526
$ws.Range("J132").Value = 4002.75
$ws.Range("K132").Value = 48265.578
$ws.Range("L132").Value = 12008.25
$ws.Range("M132").Value = -45735.578
$ws.Range("N132").Value = -17068.25
$ws.Range("H136").Value = 6573.6665
$ws.Range("I136").Value = 3052.1667
$ws.Range("J136").Value = 8921.333000000001
$ws.Range("K136").Value = 9156.500100000001
$ws.Range("L136").Value = 26763.999
$ws.Range("M136").Value = -6606.500100000001
$ws.Range("N136").Value = -31863.999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 166851.69
$ws.Range("I62").Value = 346276.34
$ws.Range("J62").Value = 5369.5
$ws.Range("K62").Value = 346276.34
$ws.Range("L62").Value = 5369.5
$ws.Range("M62").Value = -345652.34
$ws.Range("N62").Value = -6617.5
$ws.Range("H65").Value = 166851.69
$ws.Range("I65").Value = 346276.34
$ws.Range("J65").Value = 5369.5
$ws.Range("K65").Value = 1731381.7
$ws.Range("L65").Value = 26847.5
$ws.Range("M65").Value = -1728261.7
$ws.Range("N65").Value = -33087.5
$ws.Range("H107").Value = 1921.6666
$ws.Range("I107").Value = 1590.95
$ws.Range("K107").Value = 4772.85
$ws.Range("M107").Value = -2852.85
$ws.Range("H122").Value = 4304.885
$ws.Range("I122").Value = 1612.8889
$ws.Range("K122").Value = 4838.6667
$ws.Range("M122").Value = -2388.6667
$ws.Range("H126").Value = 26565.85
$ws.Range("I126").Value = 34919.785
$ws.Range("J126").Value = 7073.3335
$ws.Range("K126").Value = 104759.355
$ws.Range("L126").Value = 21220.0005
$ws.Range("M126").Value = -102289.355
$ws.Range("N126").Value = -26160.0005
$ws.Range("H132").Value = 10907.759
$ws.Range("I132").Value = 10975.187
$ws.Range("J132").Value = 10644.182
$ws.Range("K132").Value = 32925.561
$ws.Range("L132").Value = 31932.546
$ws.Range("M132").Value = -30395.561
$ws.Range("N132").Value = -36992.546
$ws.Range("H136").Value = 3788.093
$ws.Range("I136").Value = 2946.0286
$ws.Range("J136").Value = 7472.125
$ws.Range("K136").Value = 8838.085800000001
$ws.Range("L136").Value = 22416.375
$ws.Range("M136").Value = -6288.085800000001
$ws.Range("N136").Value = -27516.375
